# Auto-generated edit script for ResonatorArray.xlsx update
# Applies the small numeric tweaks (rows 15, 89) and the row-90..105
# re-sequencing / new GP+GNDfeed_bondpad row described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 15
$ws.Range("V15").Value = 0.7

# Row 89
$ws.Range("M89").Value = -0.5629999999999999
$ws.Range("O89").Value = -0.7370000000000001
$ws.Range("V89").Value = 0.7

# Row 90
$ws.Range("A90").Value = '400nm_NbWiring'
$ws.Range("C90").Value = 'Cap_to_Ind_lines'
$ws.Range("D90").Value = 0
$ws.Range("E90").Value = 0
$ws.Range("F90").Value = 1.2
$ws.Range("G90").Value = 1.3
$ws.Range("H90").Value = -0.6
$ws.Range("I90").Value = 0.6
$ws.Range("J90").Value = 0.65
$ws.Range("K90").Value = -0.65
$ws.Range("M90").Value = 0.2585
$ws.Range("O90").Value = -0.2585
$ws.Range("P90").Value = 0
$ws.Range("T90").Value = 8
$ws.Range("U90").Value = 8
$ws.Range("V90").Value = 0.7
$ws.Range("X90").Value = 11
$ws.Range("Y90").Value = 11

# Row 91
$ws.Range("C91").Value = 'GP_edge_filler_hor'
$ws.Range("E91").Value = -11.25
$ws.Range("F91").Value = 16.9
$ws.Range("G91").Value = 0.9
$ws.Range("H91").Value = -8.449999999999999
$ws.Range("I91").Value = 8.449999999999999
$ws.Range("J91").Value = -10.8
$ws.Range("K91").Value = -11.7
$ws.Range("M91").Value = 0
$ws.Range("O91").Value = 0
$ws.Range("P91").Value = -11.25
$ws.Range("T91").Value = 3
$ws.Range("U91").Value = 2
$ws.Range("V91").Value = -0
$ws.Range("X91").Value = 38.9
$ws.Range("Y91").Value = 99.88

# Row 92
$ws.Range("C92").Value = 'GP_edge_filler_vert'
$ws.Range("D92").Value = 9.75
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0.9
$ws.Range("G92").Value = 16.9
$ws.Range("H92").Value = 9.300000000000001
$ws.Range("I92").Value = 10.2
$ws.Range("J92").Value = 8.449999999999999
$ws.Range("K92").Value = -8.449999999999999
$ws.Range("O92").Value = 9.75
$ws.Range("P92").Value = 0
$ws.Range("T92").Value = 1
$ws.Range("U92").Value = 3
$ws.Range("V92").Value = -49.95
$ws.Range("X92").Value = 99.88
$ws.Range("Y92").Value = 38.9

# Row 93
$ws.Range("C93").Value = 'MSfeed_bondpad'
$ws.Range("D93").Value = 3
$ws.Range("E93").Value = -0.35
$ws.Range("F93").Value = 0.8
$ws.Range("G93").Value = 0.8
$ws.Range("H93").Value = 2.6
$ws.Range("I93").Value = 3.4
$ws.Range("J93").Value = 0.05000000000000004
$ws.Range("K93").Value = -0.75
$ws.Range("O93").Value = 3
$ws.Range("P93").Value = -0.35
$ws.Range("U93").Value = 2
$ws.Range("V93").Value = 2.996
$ws.Range("X93").Value = 0
$ws.Range("Y93").Value = 32.604

# Row 94
$ws.Range("C94").Value = 'cap_to_feed'
$ws.Range("D94").Value = 2.2
$ws.Range("E94").Value = 0.65
$ws.Range("F94").Value = 0.4
$ws.Range("H94").Value = 2
$ws.Range("I94").Value = 2.4
$ws.Range("J94").Value = 1.05
$ws.Range("K94").Value = 0.25
$ws.Range("M94").Value = -0.3
$ws.Range("N94").Value = 0.756
$ws.Range("O94").Value = 2.5
$ws.Range("P94").Value = -0.106
$ws.Range("T94").Value = 8
$ws.Range("U94").Value = 8
$ws.Range("V94").Value = -0
$ws.Range("W94").Value = 0.099
$ws.Range("X94").Value = 11
$ws.Range("Y94").Value = 11

# Row 95
$ws.Range("C95").Value = 'cap_to_gnd'
$ws.Range("E95").Value = -0.25
$ws.Range("G95").Value = 0.6
$ws.Range("J95").Value = 0.04999999999999999
$ws.Range("K95").Value = -0.55
$ws.Range("N95").Value = -0.662
$ws.Range("P95").Value = 0.412

# Row 96
$ws.Range("C96").Value = 'feedline_main'
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 12
$ws.Range("F96").Value = 11.4
$ws.Range("G96").Value = 0.4
$ws.Range("H96").Value = -5.7
$ws.Range("I96").Value = 5.7
$ws.Range("J96").Value = 12.2
$ws.Range("K96").Value = 11.8
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = 12

# Row 97
$ws.Range("C97").Value = 'main_hor_feedline_to_pad'
$ws.Range("E97").Value = 10.8
$ws.Range("F97").Value = 6.4
$ws.Range("H97").Value = -3.2
$ws.Range("I97").Value = 3.2
$ws.Range("J97").Value = 11
$ws.Range("K97").Value = 10.6
$ws.Range("N97").Value = -38.5
$ws.Range("P97").Value = 49.3
$ws.Range("T97").Value = 1
$ws.Range("U97").Value = 2
$ws.Range("W97").Value = 0
$ws.Range("X97").Value = 0
$ws.Range("Y97").Value = 77

# Row 98
$ws.Range("C98").Value = 'main_vert_feedline_to_pad'
$ws.Range("D98").Value = 10.6
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0.4
$ws.Range("G98").Value = 22.4
$ws.Range("H98").Value = 10.4
$ws.Range("I98").Value = 10.8
$ws.Range("J98").Value = 11.2
$ws.Range("K98").Value = -11.2
$ws.Range("M98").Value = 2.996
$ws.Range("N98").Value = -27.496
$ws.Range("O98").Value = 7.603999999999999
$ws.Range("P98").Value = 27.496
$ws.Range("V98").Value = 2.996
$ws.Range("Y98").Value = 54.992

# Row 99
$ws.Range("C99").Value = 'vert_main_with_corners'
$ws.Range("D99").Value = -9.4
$ws.Range("G99").Value = 11.4
$ws.Range("H99").Value = -9.6
$ws.Range("I99").Value = -9.200000000000001
$ws.Range("J99").Value = 5.7
$ws.Range("K99").Value = -5.7
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("O99").Value = -9.4
$ws.Range("P99").Value = 0
$ws.Range("U99").Value = 4
$ws.Range("V99").Value = -44
$ws.Range("W99").Value = 0.959
$ws.Range("Y99").Value = 22

# Row 100
$ws.Range("U100").Value = 3
$ws.Range("V100").Value = 44

# Row 101
$ws.Range("A101").Value = 'ILD'
$ws.Range("C101").Value = 'GP_edge_opening_hor_r'
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = -12.35
$ws.Range("F101").Value = 16.9
$ws.Range("G101").Value = 0.9
$ws.Range("H101").Value = -8.449999999999999
$ws.Range("I101").Value = 8.449999999999999
$ws.Range("J101").Value = -11.9
$ws.Range("K101").Value = -12.8
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = -12.35
$ws.Range("T101").Value = 3
$ws.Range("U101").Value = 2
$ws.Range("V101").Value = -0
$ws.Range("W101").Value = 0
$ws.Range("X101").Value = 38.9
$ws.Range("Y101").Value = 99.88

# Row 102
$ws.Range("C102").Value = 'GP_edge_opening_vert_r'
$ws.Range("D102").Value = 8.65
$ws.Range("E102").Value = 0
$ws.Range("F102").Value = 0.9
$ws.Range("G102").Value = 16.9
$ws.Range("H102").Value = 8.200000000000001
$ws.Range("I102").Value = 9.1
$ws.Range("J102").Value = 8.449999999999999
$ws.Range("K102").Value = -8.449999999999999
$ws.Range("O102").Value = 8.65
$ws.Range("P102").Value = 0
$ws.Range("T102").Value = 1
$ws.Range("U102").Value = 3
$ws.Range("V102").Value = -49.95
$ws.Range("X102").Value = 99.88
$ws.Range("Y102").Value = 38.9

# Row 103
$ws.Range("C103").Value = 'Via_to_Ground'
$ws.Range("D103").Value = 1.75
$ws.Range("E103").Value = 0.4
$ws.Range("F103").Value = 0.4
$ws.Range("G103").Value = 0.4
$ws.Range("H103").Value = 1.55
$ws.Range("I103").Value = 1.95
$ws.Range("J103").Value = 0.6000000000000001
$ws.Range("K103").Value = 0.2
$ws.Range("M103").Value = -0.3
$ws.Range("N103").Value = -0.767
$ws.Range("O103").Value = 2.05
$ws.Range("P103").Value = 1.167
$ws.Range("T103").Value = 8
$ws.Range("U103").Value = 8
$ws.Range("V103").Value = -0
$ws.Range("W103").Value = 0.099
$ws.Range("X103").Value = 11
$ws.Range("Y103").Value = 11

# Row 104
$ws.Range("C104").Value = 'reso_ILD_sub'
$ws.Range("D104").Value = 6.45
$ws.Range("E104").Value = -2.15
$ws.Range("F104").Value = 3.1
$ws.Range("G104").Value = 1.8
$ws.Range("H104").Value = 4.9
$ws.Range("I104").Value = 8
$ws.Range("J104").Value = -1.25
$ws.Range("K104").Value = -3.05
$ws.Range("M104").Value = -0.5629999999999999
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 7.013
$ws.Range("P104").Value = -2.15

# Row 105
$ws.Range("A105").Value = 'GP'
$ws.Range("C105").Value = 'GNDfeed_bondpad'
$ws.Range("D105").Value = 3
$ws.Range("E105").Value = 0.65
$ws.Range("F105").Value = 0.8
$ws.Range("G105").Value = 0.8
$ws.Range("H105").Value = 2.6
$ws.Range("I105").Value = 3.4
$ws.Range("J105").Value = 1.05
$ws.Range("K105").Value = 0.25
$ws.Range("M105").Value = 0
$ws.Range("O105").Value = 3
$ws.Range("P105").Value = 0.65
$ws.Range("T105").Value = 1
$ws.Range("U105").Value = 2
$ws.Range("V105").Value = 2.996
$ws.Range("W105").Value = 0
$ws.Range("X105").Value = 0
$ws.Range("Y105").Value = 32.604
